{"js": "// Update the date line and every \"NNN\u00f7N=\" division prompt in the table\n// to the new values from the target revision. Each \"before\" string is\n// unique in the document, so a plain exact-match search-and-replace for\n// each pair is unambiguous.\nconst replacements = [\n  [\"2024-05-15 Wednesday\", \"2024-05-16 Thursday\"],\n  [\"694\u00f78=\", \"286\u00f74=\"],\n  [\"269\u00f77=\", \"828\u00f75=\"],\n  [\"537\u00f76=\", \"424\u00f74=\"],\n  [\"930\u00f78=\", \"444\u00f78=\"],\n  [\"900\u00f75=\", \"263\u00f76=\"],\n  [\"872\u00f74=\", \"230\u00f77=\"],\n  [\"672\u00f72=\", \"171\u00f79=\"],\n  [\"453\u00f73=\", \"559\u00f76=\"],\n  [\"629\u00f76=\", \"725\u00f73=\"],\n  [\"959\u00f76=\", \"291\u00f73=\"],\n  [\"640\u00f75=\", \"459\u00f77=\"],\n  [\"445\u00f74=\", \"871\u00f78=\"],\n  [\"879\u00f74=\", \"813\u00f74=\"],\n  [\"215\u00f76=\", \"905\u00f79=\"],\n  [\"503\u00f73=\", \"981\u00f72=\"],\n  [\"482\u00f78=\", \"939\u00f73=\"],\n  [\"167\u00f72=\", \"490\u00f73=\"],\n  [\"440\u00f79=\", \"589\u00f74=\"],\n  [\"191\u00f78=\", \"370\u00f73=\"],\n  [\"521\u00f76=\", \"162\u00f76=\"],\n  [\"208\u00f77=\", \"351\u00f77=\"],\n  [\"270\u00f74=\", \"666\u00f75=\"],\n  [\"109\u00f78=\", \"829\u00f78=\"],\n  [\"292\u00f77=\", \"573\u00f72=\"],\n  [\"578\u00f78=\", \"323\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"NNN\u00f7N=\" division prompt in the table\n# to the new values from the target revision. Each \"before\" string is\n# unique in the document, so an exact Find/Replace for each pair is\n# unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-05-15 Wednesday\", \"2024-05-16 Thursday\"),\n    @(\"694\u00f78=\", \"286\u00f74=\"),\n    @(\"269\u00f77=\", \"828\u00f75=\"),\n    @(\"537\u00f76=\", \"424\u00f74=\"),\n    @(\"930\u00f78=\", \"444\u00f78=\"),\n    @(\"900\u00f75=\", \"263\u00f76=\"),\n    @(\"872\u00f74=\", \"230\u00f77=\"),\n    @(\"672\u00f72=\", \"171\u00f79=\"),\n    @(\"453\u00f73=\", \"559\u00f76=\"),\n    @(\"629\u00f76=\", \"725\u00f73=\"),\n    @(\"959\u00f76=\", \"291\u00f73=\"),\n    @(\"640\u00f75=\", \"459\u00f77=\"),\n    @(\"445\u00f74=\", \"871\u00f78=\"),\n    @(\"879\u00f74=\", \"813\u00f74=\"),\n    @(\"215\u00f76=\", \"905\u00f79=\"),\n    @(\"503\u00f73=\", \"981\u00f72=\"),\n    @(\"482\u00f78=\", \"939\u00f73=\"),\n    @(\"167\u00f72=\", \"490\u00f73=\"),\n    @(\"440\u00f79=\", \"589\u00f74=\"),\n    @(\"191\u00f78=\", \"370\u00f73=\"),\n    @(\"521\u00f76=\", \"162\u00f76=\"),\n    @(\"208\u00f77=\", \"351\u00f77=\"),\n    @(\"270\u00f74=\", \"666\u00f75=\"),\n    @(\"109\u00f78=\", \"829\u00f78=\"),\n    @(\"292\u00f77=\", \"573\u00f72=\"),\n    @(\"578\u00f78=\", \"323\u00f74=\")\n)\n\nforeach ($pair in $pairs) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $before\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $after\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n"}
